$wb = $excel.ActiveWorkbook

# --- ManageOfferCodeTest: rework the offer-code row into a percentage row, clear the old percentage row ---
$wsOfferCode = $wb.Worksheets.Item("ManageOfferCodeTest")
$wsOfferCode.Range("A4").Value = "  percentage"
$wsOfferCode.Range("B4").Value = 25
$wsOfferCode.Range("A5").Value = ""
$wsOfferCode.Range("B5").Value = ""
$wsOfferCode.Range("A5:B5").Select()

# --- ManageDeliveryBoyTest: add a new test-data row (selectedMenu / Manage Delivery Boy) ---
$wsDeliveryBoy = $wb.Worksheets.Item("ManageDeliveryBoyTest")
$wsDeliveryBoy.Range("A3").Value = "selectedMenu"
$wsDeliveryBoy.Range("B3").Value = "Manage Delivery Boy"
$wsDeliveryBoy.PageSetup.Orientation = 1

# --- Re-activate the originally selected tab / selection ---
$wsDeliveryBoy.Activate()
$wsDeliveryBoy.Range("B3").Select()
